$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: headers "n", "steps", "b"
$ws.Range("C2").Value = "n"
$ws.Range("D2").Value = "steps"
$ws.Range("E2").Value = "b"

# Row 3: starting values + formula referencing the last step result
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 25
$ws.Range("E3").Formula = "=H4"

# Row 4: the "nested cycles" chain of IF formulas
$ws.Range("C4").Value = 7
$ws.Range("D4").Formula = "=IF(C3>=1,C4+D3)"
$ws.Range("E4").Formula = "=IF(C3>=2,D4+D3)"
$ws.Range("F4").Formula = "=IF(C3>=3,E4+D3)"
$ws.Range("G4").Formula = "=IF(C3>=4,F4+D3)"
$ws.Range("H4").Formula = "=IF(C3>=5,G4+D3)"

# New H4/H5 cells need the same box-border formatting as the rest of their row
# (PasteSpecial xlPasteFormats=-4122 copies style only, leaving formulas/values intact)
$ws.Range("G4").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4122) | Out-Null
$ws.Range("G5").Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Add a watch window entry on E4 (matches the cellWatches element in the file)
try {
    $excel.Watches.Add($ws.Range("E4")) | Out-Null
} catch {
}

# Move/record the active selection
$ws.Range("E8").Select() | Out-Null

Write-Output "done"
